$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materi")

# Row 2
$ws.Cells.Item(2,14).Value = "TC-Membuat Video Pembahasan-001"

# Row 3
$ws.Cells.Item(3,14).Value = "TC-Membuat Video Pembahasan-002"

# Row 4
$ws.Cells.Item(4,3).Value = "java-logo.jpg"
$ws.Cells.Item(4,4).Value = "padepokan79oke.com"
$ws.Cells.Item(4,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(4,6).Value = "Fundamental Java"
$ws.Cells.Item(4,7).Value = "Perulangan"
$ws.Cells.Item(4,9).Value = "Foreach"
$ws.Cells.Item(4,10).Value = "Publik"
$ws.Cells.Item(4,11).Value = "JAVA"
$ws.Cells.Item(4,12).Value = "'99"
$ws.Cells.Item(4,14).Value = "TC-Membuat Video Pembahasan-003"

# Row 5
$ws.Cells.Item(5,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(5,4).Value = "padepokan79oke.com"
$ws.Cells.Item(5,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(5,6).Value = "Fundamental Java"
$ws.Cells.Item(5,7).Value = "Perulangan"
$ws.Cells.Item(5,9).Value = "Foreach"
$ws.Cells.Item(5,10).Value = "Publik"
$ws.Cells.Item(5,11).Value = "JAVA"
$ws.Cells.Item(5,12).Value = "'99"
$ws.Cells.Item(5,14).Value = "TC-Membuat Video Pembahasan-004"

# Row 6
$ws.Cells.Item(6,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(6,3).Value = "ukuran tidak sesuai"
$ws.Cells.Item(6,4).Value = "padepokan79oke.com"
$ws.Cells.Item(6,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(6,6).Value = "Fundamental Java"
$ws.Cells.Item(6,7).Value = "Perulangan"
$ws.Cells.Item(6,9).Value = "Foreach"
$ws.Cells.Item(6,10).Value = "Publik"
$ws.Cells.Item(6,11).Value = "JAVA"
$ws.Cells.Item(6,12).Value = "'99"
$ws.Cells.Item(6,14).Value = "TC-Membuat Video Pembahasan-005"

# Row 7
$ws.Cells.Item(7,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(7,3).Value = "java-logo.jpg"
$ws.Cells.Item(7,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(7,6).Value = "Fundamental Java"
$ws.Cells.Item(7,7).Value = "Perulangan"
$ws.Cells.Item(7,9).Value = "Foreach"
$ws.Cells.Item(7,10).Value = "Publik"
$ws.Cells.Item(7,11).Value = "JAVA"
$ws.Cells.Item(7,12).Value = "'99"
$ws.Cells.Item(7,14).Value = "TC-Membuat Video Pembahasan-006"

# Row 8
$ws.Cells.Item(8,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(8,3).Value = "java-logo.jpg"
$ws.Cells.Item(8,4).Value = "padepokan79oke.com"
$ws.Cells.Item(8,6).Value = "Fundamental Java"
$ws.Cells.Item(8,7).Value = "Perulangan"
$ws.Cells.Item(8,9).Value = "Foreach"
$ws.Cells.Item(8,10).Value = "Publik"
$ws.Cells.Item(8,11).Value = "JAVA"
$ws.Cells.Item(8,12).Value = "'99"
$ws.Cells.Item(8,14).Value = "TC-Membuat Video Pembahasan-007"

# Row 9
$ws.Cells.Item(9,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(9,3).Value = "java-logo.jpg"
$ws.Cells.Item(9,4).Value = "padepokan79oke.com"
$ws.Cells.Item(9,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(9,10).Value = "Publik"
$ws.Cells.Item(9,11).Value = "JAVA"
$ws.Cells.Item(9,12).Value = "'99"
$ws.Cells.Item(9,14).Value = "TC-Membuat Video Pembahasan-008"

# Row 10
$ws.Cells.Item(10,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(10,3).Value = "java-logo.jpg"
$ws.Cells.Item(10,4).Value = "padepokan79oke.com"
$ws.Cells.Item(10,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(10,6).Value = "Fundamental Java"
$ws.Cells.Item(10,10).Value = "Publik"
$ws.Cells.Item(10,11).Value = "JAVA"
$ws.Cells.Item(10,12).Value = "'99"
$ws.Cells.Item(10,14).Value = "TC-Membuat Video Pembahasan-009"

# Row 11
$ws.Cells.Item(11,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(11,3).Value = "java-logo.jpg"
$ws.Cells.Item(11,4).Value = "padepokan79oke.com"
$ws.Cells.Item(11,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(11,6).Value = "Fundamental Java"
$ws.Cells.Item(11,7).Value = "Perulangan"
$ws.Cells.Item(11,10).Value = "Publik"
$ws.Cells.Item(11,11).Value = "JAVA"
$ws.Cells.Item(11,12).Value = "'99"
$ws.Cells.Item(11,14).Value = "TC-Membuat Video Pembahasan-010"

# Row 12
$ws.Cells.Item(12,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(12,3).Value = "java-logo.jpg"
$ws.Cells.Item(12,4).Value = "padepokan79oke.com"
$ws.Cells.Item(12,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(12,6).Value = "Fundamental Java"
$ws.Cells.Item(12,7).Value = "Perulangan"
$ws.Cells.Item(12,9).Value = "Foreach"
$ws.Cells.Item(12,10).Value = "Publik"
$ws.Cells.Item(12,12).Value = "'99"
$ws.Cells.Item(12,14).Value = "TC-Membuat Video Pembahasan-011"

# Row 13
$ws.Cells.Item(13,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(13,3).Value = "java-logo.jpg"
$ws.Cells.Item(13,4).Value = "padepokan79oke.com"
$ws.Cells.Item(13,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(13,6).Value = "Fundamental Java"
$ws.Cells.Item(13,7).Value = "Perulangan"
$ws.Cells.Item(13,9).Value = "Foreach"
$ws.Cells.Item(13,10).Value = "Publik"
$ws.Cells.Item(13,11).Value = "JAVA"
$ws.Cells.Item(13,14).Value = "TC-Membuat Video Pembahasan-012"

# Row 14
$ws.Cells.Item(14,14).Value = "TC-Membuat Video Pembahasan-013"

# Row 15
$ws.Cells.Item(15,2).Value = "Materi pembahasan Katalon"
$ws.Cells.Item(15,3).Value = "java-logo.jpg"
$ws.Cells.Item(15,4).Value = "padepokan79oke.com"
$ws.Cells.Item(15,5).Value = "Fundamental Katalon yang mengenalkan java dari dasar."
$ws.Cells.Item(15,6).Value = "Fundamental Java"
$ws.Cells.Item(15,7).Value = "Perulangan"
$ws.Cells.Item(15,9).Value = "Foreach"
$ws.Cells.Item(15,10).Value = "Publik"
$ws.Cells.Item(15,11).Value = "JAVA"
$ws.Cells.Item(15,12).Value = "'99"
$ws.Cells.Item(15,14).Value = "TC-Membuat Video Pembahasan-014"

# Row 16
$ws.Cells.Item(16,1).Value = "'17"
$ws.Cells.Item(16,2).Value = "Materi pembahasan Selenium"
$ws.Cells.Item(16,3).Value = "java-logo.jpg"
$ws.Cells.Item(16,4).Value = "padepokan79oke2.com"
$ws.Cells.Item(16,5).Value = "Fundamental Selenium yang mengenalkan java dari dasar."
$ws.Cells.Item(16,6).Value = "Fundamental Java"
$ws.Cells.Item(16,7).Value = "Perulangan"
$ws.Cells.Item(16,9).Value = "Foreach"
$ws.Cells.Item(16,10).Value = "Sembunyi"
$ws.Cells.Item(16,11).Value = "FUNDAMENTAL"
$ws.Cells.Item(16,12).Value = "'89"
$ws.Cells.Item(16,14).Value = "TC-Mengubah Video Pembahasan-001"

# Row 17
$ws.Cells.Item(17,1).Value = "'17"
$ws.Cells.Item(17,14).Value = "TC-Mengubah Video Pembahasan-002"

# Row 18
$ws.Cells.Item(18,1).Value = "'17"
$ws.Cells.Item(18,3).Value = "java-logo.jpg"
$ws.Cells.Item(18,4).Value = "padepokan79oke2.com"
$ws.Cells.Item(18,5).Value = "Fundamental Selenium yang mengenalkan java dari dasar."
$ws.Cells.Item(18,6).Value = "Fundamental Java"
$ws.Cells.Item(18,7).Value = "Perulangan"
$ws.Cells.Item(18,9).Value = "Foreach"
$ws.Cells.Item(18,10).Value = "Sembunyi"
$ws.Cells.Item(18,11).Value = "FUNDAMENTAL"
$ws.Cells.Item(18,12).Value = "'89"
$ws.Cells.Item(18,14).Value = "TC-Mengubah Video Pembahasan-003"

# Row 19
$ws.Cells.Item(19,1).Value = "'17"
$ws.Cells.Item(19,2).Value = "Materi pembahasan Selenium"
$ws.Cells.Item(19,3).Value = "tidak sesuai"
$ws.Cells.Item(19,4).Value = "padepokan79oke2.com"
$ws.Cells.Item(19,5).Value = "Fundamental Selenium yang mengenalkan java dari dasar."
$ws.Cells.Item(19,6).Value = "Fundamental Java"
$ws.Cells.Item(19,7).Value = "Perulangan"
$ws.Cells.Item(19,9).Value = "Foreach"
$ws.Cells.Item(19,10).Value = "Sembunyi"
$ws.Cells.Item(19,11).Value = "FUNDAMENTAL"
$ws.Cells.Item(19,12).Value = "'89"
$ws.Cells.Item(19,14).Value = "TC-Mengubah Video Pembahasan-004"

# Row 20
$ws.Cells.Item(20,1).Value = "'17"
$ws.Cells.Item(20,2).Value = "Materi pembahasan Selenium"
$ws.Cells.Item(20,3).Value = "java-logo.jpg"
$ws.Cells.Item(20,5).Value = "Fundamental Selenium yang mengenalkan java dari dasar."
$ws.Cells.Item(20,6).Value = "Fundamental Java"
$ws.Cells.Item(20,7).Value = "Perulangan"
$ws.Cells.Item(20,9).Value = "Foreach"
$ws.Cells.Item(20,10).Value = "Sembunyi"
$ws.Cells.Item(20,11).Value = "FUNDAMENTAL"
$ws.Cells.Item(20,12).Value = "'89"
$ws.Cells.Item(20,14).Value = "TC-Mengubah Video Pembahasan-005"

# Row 21
$ws.Cells.Item(21,1).Value = "'17"
$ws.Cells.Item(21,2).Value = "Materi pembahasan Selenium"
$ws.Cells.Item(21,3).Value = "java-logo.jpg"
$ws.Cells.Item(21,4).Value = "padepokan79oke2.com"
$ws.Cells.Item(21,6).Value = "Fundamental Java"
$ws.Cells.Item(21,7).Value = "Perulangan"
$ws.Cells.Item(21,9).Value = "Foreach"
$ws.Cells.Item(21,10).Value = "Sembunyi"
$ws.Cells.Item(21,11).Value = "FUNDAMENTAL"
$ws.Cells.Item(21,12).Value = "'89"
$ws.Cells.Item(21,14).Value = "TC-Mengubah Video Pembahasan-006"

# Row 22
$ws.Cells.Item(22,1).Value = "'17"
$ws.Cells.Item(22,2).Value = "Materi pembahasan Selenium"
$ws.Cells.Item(22,3).Value = "java-logo.jpg"
$ws.Cells.Item(22,4).Value = "padepokan79oke2.com"
$ws.Cells.Item(22,5).Value = "Fundamental Selenium yang mengenalkan java dari dasar."
$ws.Cells.Item(22,6).Value = "Fundamental Java"
$ws.Cells.Item(22,7).Value = "Perulangan"
$ws.Cells.Item(22,9).Value = "Foreach"
$ws.Cells.Item(22,10).Value = "Sembunyi"
$ws.Cells.Item(22,11).Value = "FUNDAMENTAL"
$ws.Cells.Item(22,14).Value = "TC-Mengubah Video Pembahasan-007"

# Row 23
$ws.Cells.Item(23,1).Value = "'17"
$ws.Cells.Item(23,2).Value = "Materi pembahasan Selenium"
$ws.Cells.Item(23,3).Value = "java-logo.jpg"
$ws.Cells.Item(23,4).Value = "padepokan79oke2.com"
$ws.Cells.Item(23,5).Value = "Fundamental Selenium yang mengenalkan java dari dasar."
$ws.Cells.Item(23,6).Value = "Fundamental Java"
$ws.Cells.Item(23,7).Value = "Perulangan"
$ws.Cells.Item(23,9).Value = "Foreach"
$ws.Cells.Item(23,10).Value = "Sembunyi"
$ws.Cells.Item(23,12).Value = "'89"
$ws.Cells.Item(23,14).Value = "TC-Mengubah Video Pembahasan-008"

# Row 24
$ws.Cells.Item(24,1).Value = "'17"
$ws.Cells.Item(24,2).Value = "Materi pembahasan Selenium"
$ws.Cells.Item(24,3).Value = "java-logo.jpg"
$ws.Cells.Item(24,4).Value = "padepokan79oke2.com"
$ws.Cells.Item(24,5).Value = "Fundamental Selenium yang mengenalkan java dari dasar."
$ws.Cells.Item(24,6).Value = "Fundamental Java"
$ws.Cells.Item(24,7).Value = "Perulangan"
$ws.Cells.Item(24,9).Value = "Foreach"
$ws.Cells.Item(24,10).Value = "Sembunyi"
$ws.Cells.Item(24,11).Value = "FUNDAMENTAL"
$ws.Cells.Item(24,12).Value = "'89"
$ws.Cells.Item(24,14).Value = "TC-Mengubah Video Pembahasan-009"

# Row 25
$ws.Cells.Item(25,1).Value = "'17"
$ws.Cells.Item(25,14).Value = "TC-Menghapus Video Pembahasan-001"

# Row 26
$ws.Cells.Item(26,1).Value = "'17"
$ws.Cells.Item(26,14).Value = "TC-Menghapus Video Pembahasan-002"

# Row 27
$ws.Cells.Item(27,1).Value = "'17"
$ws.Cells.Item(27,14).Value = "TC-Menghapus Video Pembahasan-003"

# Row 28
$ws.Cells.Item(28,8).Value = "Pembahasan"
$ws.Cells.Item(28,14).Value = "TC-Mencari Video Pembahasan-001"

# Row 29
$ws.Cells.Item(29,8).Value = "Pendahuluan"
$ws.Cells.Item(29,14).Value = "TC-Mencari Video Pembahasan-002"

# Row 30
$ws.Cells.Item(30,8).Value = "Java"
$ws.Cells.Item(30,14).Value = "TC-Mencari Video Pembahasan-003"

# Row 31
$ws.Cells.Item(31,8).Value = "Foreeach"
$ws.Cells.Item(31,14).Value = "TC-Mencari Video Pembahasan-004"

# Row 32
$ws.Cells.Item(32,8).Value = "JAVA"
$ws.Cells.Item(32,14).Value = "TC-Mencari Video Pembahasan-005"

# Row 33
$ws.Cells.Item(33,8).Value = "Sembunyi"
$ws.Cells.Item(33,14).Value = "TC-Mencari Video Pembahasan-006"

# Row 34
$ws.Cells.Item(34,8).Value = "Jawa"
$ws.Cells.Item(34,14).Value = "TC-Mencari Video Pembahasan-007"

# Row 35
$ws.Cells.Item(35,8).Value = "Pembahasan"
$ws.Cells.Item(35,14).Value = "TC-Mencari Video Pembahasan-008"

# Row 36
$ws.Cells.Item(36,8).Value = "Sembunyi"
$ws.Cells.Item(36,14).Value = "TC-Mencari Video Pembahasan-009"

# Row 37
$ws.Cells.Item(37,8).Value = "Foreeach"
$ws.Cells.Item(37,14).Value = "TC-Mencari Video Pembahasan-010"

# Row 38
$ws.Cells.Item(38,8).Value = "Sembunyi"
$ws.Cells.Item(38,14).Value = "TC-Mencari Video Pembahasan-011"

# Row 39
$ws.Cells.Item(39,8).Value = "Java"
$ws.Cells.Item(39,14).Value = "TC-Mencari Video Pembahasan-012"

# Row 40
$ws.Cells.Item(40,8).Value = "Sembunyi"
$ws.Cells.Item(40,14).Value = "TC-Mencari Video Pembahasan-013"

# Row 41
$ws.Cells.Item(41,8).Value = "Pendahuluan"
$ws.Cells.Item(41,14).Value = "TC-Mencari Video Pembahasan-014"

# Row 42
$ws.Cells.Item(42,8).Value = "Sembunyi"
$ws.Cells.Item(42,14).Value = "TC-Mencari Video Pembahasan-015"

# Row 43
$ws.Cells.Item(43,8).Value = "JAVA"
$ws.Cells.Item(43,14).Value = "TC-Mencari Video Pembahasan-016"

# Row 44
$ws.Cells.Item(44,8).Value = "Sembunyi"
$ws.Cells.Item(44,14).Value = "TC-Mencari Video Pembahasan-017"

# Row 45
$ws.Cells.Item(45,8).Value = "Sembunyi"
$ws.Cells.Item(45,14).Value = "TC-Mencari Video Pembahasan-018"

# Row 46
$ws.Cells.Item(46,8).Value = "Java"
$ws.Cells.Item(46,14).Value = "TC-Mencari Video Pembahasan-019"

# Row 13: L13 keeps quote-prefix style with no value
$ws.Cells.Item(13,12).Value = "'99"
$ws.Cells.Item(13,12).Value = ""

# Row 22: L22 keeps quote-prefix style with no value
$ws.Cells.Item(22,12).Value = "'99"
$ws.Cells.Item(22,12).Value = ""
